# Scheduled market-data refresh: update currentAveragePrice / LevePrice /
# LeveProfit columns (H:N) for the affected Levequest rows across each
# job sheet (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 40.5
$ws.Range("I11").Value = 40.5
$ws.Range("K11").Value = 40.5
$ws.Range("M11").Value = 99.5
$ws.Range("H40").Value = 10085.833
$ws.Range("J40").Value = 11557.7
$ws.Range("L40").Value = 11557.7
$ws.Range("N40").Value = -11907.7
$ws.Range("H43").Value = 5131.8
$ws.Range("I43").Value = 3914.75
$ws.Range("J43").Value = 10000
$ws.Range("K43").Value = 3914.75
$ws.Range("L43").Value = 10000
$ws.Range("M43").Value = -3845.75
$ws.Range("N43").Value = -10138
$ws.Range("H69").Value = 16680.143
$ws.Range("J69").Value = 5582.6665
$ws.Range("L69").Value = 16747.9995
$ws.Range("N69").Value = -18495.9995
$ws.Range("H72").Value = 16680.143
$ws.Range("J72").Value = 5582.6665
$ws.Range("L72").Value = 50243.9985
$ws.Range("N72").Value = -58979.9985
$ws.Range("H74").Value = 3899.7334
$ws.Range("I74").Value = 3465.9167
$ws.Range("K74").Value = 3465.9167
$ws.Range("M74").Value = -2529.9167
$ws.Range("H77").Value = 3899.7334
$ws.Range("I77").Value = 3465.9167
$ws.Range("K77").Value = 17329.5835
$ws.Range("M77").Value = -12649.5835
$ws.Range("H96").Value = 125539.375
$ws.Range("I96").Value = 143202.14
$ws.Range("K96").Value = 429606.42
$ws.Range("M96").Value = -428233.42
$ws.Range("H100").Value = 3155.4
$ws.Range("I100").Value = 0
$ws.Range("J100").Value = 3155.4
$ws.Range("K100").Value = 0
$ws.Range("L100").ClearContents()
$ws.Range("M100").ClearContents()
$ws.Range("N100").Value = -4237.4
$ws.Range("H116").Value = 1674933.9
$ws.Range("J116").Value = 8342982.5
$ws.Range("L116").Value = 8342982.5
$ws.Range("N116").Value = -8349866.5
$ws.Range("H135").Value = 1370.8
$ws.Range("I135").Value = 1292.2778
$ws.Range("K135").Value = 11630.5002
$ws.Range("M135").Value = -9095.5002
$ws.Range("H138").Value = 2424.5588
$ws.Range("J138").Value = 2734.9412
$ws.Range("L138").Value = 8204.8236
$ws.Range("N138").Value = -18484.8236

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 9050.823
$ws.Range("I45").Value = 10766.615
$ws.Range("J45").Value = 3474.5
$ws.Range("K45").Value = 10766.615
$ws.Range("L45").Value = 3474.5
$ws.Range("M45").Value = -10389.615
$ws.Range("N45").Value = -4228.5
$ws.Range("H62").Value = 1249
$ws.Range("J62").Value = 1249
$ws.Range("L62").Value = 1249
$ws.Range("N62").Value = -2497
$ws.Range("H65").Value = 1249
$ws.Range("J65").Value = 1249
$ws.Range("L65").Value = 3747
$ws.Range("N65").Value = -9987
$ws.Range("H76").Value = 129995.4
$ws.Range("J76").Value = 129995.4
$ws.Range("L76").Value = 129995.4
$ws.Range("N76").Value = -130671.4
$ws.Range("H79").Value = 129995.4
$ws.Range("J79").Value = 129995.4
$ws.Range("L79").Value = 129995.4
$ws.Range("N79").Value = -132335.4
$ws.Range("H102").Value = 135246.56
$ws.Range("I102").Value = 159136.86
$ws.Range("K102").Value = 159136.86
$ws.Range("M102").Value = -157514.86
$ws.Range("H122").Value = 7765.4
$ws.Range("I122").Value = 7765.4
$ws.Range("K122").Value = 23296.2
$ws.Range("M122").Value = -20846.2

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2314.9565
$ws.Range("I20").Value = 2016
$ws.Range("K20").Value = 2016
$ws.Range("M20").Value = -1769
$ws.Range("H105").Value = 86395
$ws.Range("I105").Value = 144785.72
$ws.Range("K105").Value = 144785.72
$ws.Range("M105").Value = -143038.72
$ws.Range("H134").Value = 2720.2
$ws.Range("I134").Value = 1695.3334
$ws.Range("J134").Value = 3666.2307
$ws.Range("K134").Value = 5086.0002
$ws.Range("L134").Value = 10998.6921
$ws.Range("M134").Value = -2551.0002
$ws.Range("N134").Value = -16068.6921

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H94").Value = 6324.524
$ws.Range("I94").Value = 13245.444
$ws.Range("J94").Value = 1133.8334
$ws.Range("K94").Value = 13245.444
$ws.Range("L94").Value = 1133.8334
$ws.Range("M94").Value = -12794.444
$ws.Range("N94").Value = -2035.8334
$ws.Range("H132").Value = 2059.1052
$ws.Range("I132").Value = 2064
$ws.Range("J132").Value = 2040.75
$ws.Range("K132").Value = 6192
$ws.Range("L132").Value = 6122.25
$ws.Range("M132").Value = -3662
$ws.Range("N132").Value = -11182.25
$ws.Range("H134").Value = 65333
$ws.Range("I134").Value = 3157.7
$ws.Range("J134").Value = 168958.5
$ws.Range("K134").Value = 9473.099999999999
$ws.Range("L134").Value = 506875.5
$ws.Range("M134").Value = -6938.099999999999
$ws.Range("N134").Value = -511945.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 4284.2144
$ws.Range("I3").Value = 1197.6666
$ws.Range("J3").Value = 9840
$ws.Range("K3").Value = 3592.9998
$ws.Range("L3").Value = 29520
$ws.Range("M3").Value = -3480.9998
$ws.Range("N3").Value = -29744
$ws.Range("H4").Value = 112825310
$ws.Range("I4").Value = 131514530
$ws.Range("K4").Value = 394543590
$ws.Range("M4").Value = -394543478
$ws.Range("H10").Value = 246.71428
$ws.Range("J10").Value = 1000
$ws.Range("L10").Value = 3000
$ws.Range("N10").Value = -3278
$ws.Range("H54").Value = 3996.6667
$ws.Range("I54").Value = 0
$ws.Range("J54").Value = 3996.6667
$ws.Range("K54").Value = 0
$ws.Range("L54").ClearContents()
$ws.Range("M54").ClearContents()
$ws.Range("N54").Value = -13108.0001
$ws.Range("H58").Value = 1290
$ws.Range("I58").Value = 632.5
$ws.Range("J58").Value = 1947.5
$ws.Range("K58").Value = 1897.5
$ws.Range("L58").Value = 5842.5
$ws.Range("M58").Value = -1769.5
$ws.Range("N58").Value = -6098.5
$ws.Range("H87").Value = 100
$ws.Range("I87").Value = 100
$ws.Range("J87").Value = 0
$ws.Range("K87").Value = 300
$ws.Range("L87").ClearContents()
$ws.Range("N87").ClearContents()
$ws.Range("M87").Value = 948
$ws.Range("H90").Value = 100
$ws.Range("I90").Value = 100
$ws.Range("J90").Value = 0
$ws.Range("K90").Value = 900
$ws.Range("L90").ClearContents()
$ws.Range("N90").ClearContents()
$ws.Range("M90").Value = 5340
$ws.Range("H116").Value = 1500
$ws.Range("I116").Value = 1500
$ws.Range("K116").Value = 4500
$ws.Range("M116").Value = -1058
$ws.Range("H121").Value = 49973.285
$ws.Range("J121").Value = 55086.26
$ws.Range("L121").Value = 165258.78
$ws.Range("N121").Value = -167878.78
$ws.Range("H124").Value = 0
$ws.Range("I124").Value = 0
$ws.Range("K124").Value = 0
$ws.Range("M124").ClearContents()
$ws.Range("H125").Value = 19999
$ws.Range("J125").Value = 0
$ws.Range("L125").Value = 0
$ws.Range("N125").ClearContents()
$ws.Range("H127").Value = 2445.8333
$ws.Range("J127").Value = 2445.8333
$ws.Range("L127").Value = 7337.499899999999
$ws.Range("N127").Value = -17257.4999
$ws.Range("H139").Value = 15207.368
$ws.Range("I139").Value = 6250
$ws.Range("K139").Value = 18750
$ws.Range("M139").Value = -13610

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H53").Value = 70000000
$ws.Range("I53").Value = 70000000
$ws.Range("J53").Value = 0
$ws.Range("K53").Value = 70000000
$ws.Range("L53").ClearContents()
$ws.Range("N53").ClearContents()
$ws.Range("M53").Value = -69999369
$ws.Range("H122").Value = 16224.9
$ws.Range("J122").Value = 5784
$ws.Range("L122").Value = 17352
$ws.Range("N122").Value = -22252
$ws.Range("H126").Value = 5766.6665
$ws.Range("I126").Value = 0
$ws.Range("K126").Value = 0
$ws.Range("M126").ClearContents()
$ws.Range("H132").Value = 4783.0312
$ws.Range("I132").Value = 3883.5652
$ws.Range("K132").Value = 11650.6956
$ws.Range("M132").Value = -9120.695599999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 6948814
$ws.Range("I40").Value = 3992.3333
$ws.Range("J40").Value = 27783278
$ws.Range("K40").Value = 3992.3333
$ws.Range("L40").Value = 27783278
$ws.Range("M40").Value = -3856.3333
$ws.Range("N40").Value = -27783550
$ws.Range("H50").Value = 6000
$ws.Range("I50").Value = 5000
$ws.Range("J50").Value = 8000
$ws.Range("K50").Value = 5000
$ws.Range("L50").Value = 8000
$ws.Range("M50").Value = -4363
$ws.Range("N50").Value = -9274
$ws.Range("H54").Value = 35000
$ws.Range("J54").Value = 35000
$ws.Range("L54").Value = 35000
$ws.Range("N54").Value = -36288
$ws.Range("H82").Value = 1113.5264
$ws.Range("I82").Value = 1131.75
$ws.Range("J82").Value = 1016.3333
$ws.Range("K82").Value = 1131.75
$ws.Range("L82").Value = 1016.3333
$ws.Range("M82").Value = -770.75
$ws.Range("N82").Value = -1738.3333
$ws.Range("H85").Value = 1113.5264
$ws.Range("I85").Value = 1131.75
$ws.Range("J85").Value = 1016.3333
$ws.Range("K85").Value = 1131.75
$ws.Range("L85").Value = 1016.3333
$ws.Range("M85").Value = 116.25
$ws.Range("N85").Value = -3512.3333
$ws.Range("H136").Value = 6340.346
$ws.Range("J136").Value = 6379.6924
$ws.Range("L136").Value = 19139.0772
$ws.Range("N136").Value = -24239.0772

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H40").Value = 2999.6667
$ws.Range("I40").Value = 2999.5
$ws.Range("K40").Value = 2999.5
$ws.Range("M40").Value = -2850.5
$ws.Range("H132").Value = 1403931.8
$ws.Range("I132").Value = 1155.125
$ws.Range("J132").Value = 6213451.5
$ws.Range("K132").Value = 3465.375
$ws.Range("L132").Value = 18640354.5
$ws.Range("M132").Value = -935.375
$ws.Range("N132").Value = -18645414.5